# "Generate Report for Archive"
# The localization report was re-generated: the shared status string moved
# from "Ready for handoff" to "In Translation" everywhere it is used
# (Overview!E2/F2 and the per-language Status column on the "zh-cn" and
# "de-de" sheets), and the Status-related columns were re-autofit to the
# new (shorter) text, so they got narrower.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# Re-fit the columns that held the status text so they shrink to match the
# new, shorter value.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5
